$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the "last changed" date for every data row.
# Update it from 2023-09-01 (serial 45170) to 2023-09-05 (serial 45174)
# for every data row (rows 2 through 158).
for ($row = 2; $row -le 158; $row++) {
    $ws.Cells.Item($row, 3).Value = 45174
}
